# Weekly fruit/vegetable data refresh for "Fruta, Terminal Hortofrutícola
# Agro Chillán - Plátano": three new price observations (dated 45147) are
# inserted at the top of the data block (row 970), pushing the existing
# rows down by three (970-1024 -> 973-1027).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right above the current row 970, shifting all
# subsequent rows (and the sheet dimension) down.
$ws.Rows("970:972").Insert()

# --- New row 970 --------------------------------------------------------
$ws.Range("A970").Value = 7
$ws.Range("B970").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C970").Value = "Ñuble"
$ws.Range("D970").Value = 45147
$ws.Range("E970").Value = 16
$ws.Range("F970").Value = "Fruta"
$ws.Range("G970").Value = 100108
$ws.Range("H970").Value = "Tropicales y subtropicales"
$ws.Range("I970").Value = 100108006
$ws.Range("J970").Value = "Plátano"
$ws.Range("K970").Value = "Sin especificar"
$ws.Range("L970").Value = "Maduro"
$ws.Range("M970").Value = 300
$ws.Range("N970").Value = 15000
$ws.Range("O970").Value = 15000
$ws.Range("P970").Value = 15000
$ws.Range("Q970").Value = "$/caja 20 kilos"
$ws.Range("R970").Value = "Ecuador"
$ws.Range("S970").Value = 750
$ws.Range("T970").Value = 20

# --- New row 971 --------------------------------------------------------
$ws.Range("A971").Value = 7
$ws.Range("B971").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C971").Value = "Ñuble"
$ws.Range("D971").Value = 45147
$ws.Range("E971").Value = 16
$ws.Range("F971").Value = "Fruta"
$ws.Range("G971").Value = 100108
$ws.Range("H971").Value = "Tropicales y subtropicales"
$ws.Range("I971").Value = 100108006
$ws.Range("J971").Value = "Plátano"
$ws.Range("K971").Value = "Sin especificar"
$ws.Range("L971").Value = "Primera"
$ws.Range("M971").Value = 250
$ws.Range("N971").Value = 16000
$ws.Range("O971").Value = 16000
$ws.Range("P971").Value = 16000
$ws.Range("Q971").Value = "$/caja 20 kilos"
$ws.Range("R971").Value = "Ecuador"
$ws.Range("S971").Value = 800
$ws.Range("T971").Value = 20

# --- New row 972 --------------------------------------------------------
$ws.Range("A972").Value = 7
$ws.Range("B972").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C972").Value = "Ñuble"
$ws.Range("D972").Value = 45147
$ws.Range("E972").Value = 16
$ws.Range("F972").Value = "Fruta"
$ws.Range("G972").Value = 100108
$ws.Range("H972").Value = "Tropicales y subtropicales"
$ws.Range("I972").Value = 100108006
$ws.Range("J972").Value = "Plátano"
$ws.Range("K972").Value = "Sin especificar"
$ws.Range("L972").Value = "Primera Pintón"
$ws.Range("M972").Value = 300
$ws.Range("N972").Value = 17000
$ws.Range("O972").Value = 17000
$ws.Range("P972").Value = 17000
$ws.Range("Q972").Value = "$/caja 20 kilos"
$ws.Range("R972").Value = "Ecuador"
$ws.Range("S972").Value = 850
$ws.Range("T972").Value = 20
